$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Vega Modelo de Temuco" (Piña,
# Caramelo, Primera) on 2021-10-07 (serial 44476). It belongs right after
# the existing row 246, so insert a fresh row at position 247 — this shifts
# the previous rows 247:286 down to 248:287, preserving all of their data
# and formatting (incl. the date style on column D).
$ws.Rows.Item(247).Insert()

# Fill in the newly inserted row with the new record's values.
$ws.Range("A247").Value = 10
$ws.Range("B247").Value = "Vega Modelo de Temuco"
$ws.Range("C247").Value = "La Araucanía"
$ws.Range("D247").Value = 44476
$ws.Range("E247").Value = 9
$ws.Range("F247").Value = "Fruta"
$ws.Range("G247").Value = 100108
$ws.Range("H247").Value = "Tropicales y subtropicales"
$ws.Range("I247").Value = 100108005
$ws.Range("J247").Value = "Piña"
$ws.Range("K247").Value = "Caramelo"
$ws.Range("L247").Value = "Primera"
$ws.Range("M247").Value = 200
$ws.Range("N247").Value = 20000
$ws.Range("O247").Value = 20000
$ws.Range("P247").Value = 20000
$ws.Range("Q247").Value = "$/caja 12 unidades"
$ws.Range("R247").Value = "Ecuador"
$ws.Range("S247").Value = 1667
$ws.Range("T247").Value = 12
